# Swap the B/D/E/F/G values between each pair of adjacent rows listed below.
# (A and C columns — serial number and item name — stay put; only the
# per-batch stock figures swap rows.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(127, 128),
    @(227, 228),
    @(229, 230),
    @(243, 244),
    @(322, 323),
    @(366, 367),
    @(380, 381),
    @(382, 383),
    @(385, 386),
    @(442, 443),
    @(572, 573)
)

$cols = @("B", "D", "E", "F", "G")

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"

        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2

        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}
